# Update NATMI TPM-derived values in the "Tnc-Itga7" data sheet.
# Only the data cells in columns G:T for rows 2-10 change; all other
# cells (identifiers in A:F, K:L, headers, styles) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, taken directly from the target OOXML.
$updates = @{
    2 = @{ G=0.1347866666666667; H=0.40436; I=0.03419045085634245; J=0.03419045085634244;
           M=2.750415333333333; N=8.251245999999998; O=0.04811444325525444; P=0.04811444325525444;
           Q=0.3707193147288888; R=3.336473832559999; S=0.001645054507599054; T=0.001645054507599054 }
    3 = @{ G=0.1347866666666667; H=0.40436; I=0.03419045085634245; J=0.03419045085634244;
           O=0.01415294505639593; P=0.01415294505639593;
           Q=0.1090477149422222; R=0.9814294344800001; S=0.0004838955724232197; T=0.0004838955724232197 }
    4 = @{ G=0.1347866666666667; H=0.40436; I=0.03419045085634245; J=0.03419045085634244;
           O=0.9377326116883496; P=0.9377326116883496;
           Q=7.225181622902221; R=65.02663460612; S=0.03206150077632017; T=0.03206150077632017 }
    5 = @{ I=0.3318597741685039; J=0.3318597741685039;
           M=2.750415333333333; N=8.251245999999998; O=0.04811444325525444; P=0.04811444325525444;
           Q=3.598280367309333; R=32.38452330578399; S=0.01596724827293204; T=0.01596724827293203 }
    6 = @{ I=0.3318597741685039; J=0.3318597741685039;
           O=0.01415294505639593; P=0.01415294505639593;
           S=0.004696793150234797; T=0.004696793150234796 }
    7 = @{ I=0.3318597741685039; J=0.3318597741685039;
           O=0.9377326116883496; P=0.9377326116883496;
           S=0.3111957327453371; T=0.3111957327453371 }
    8 = @{ I=0.6339497749751537; J=0.6339497749751537;
           M=2.750415333333333; N=8.251245999999998; O=0.04811444325525444; P=0.04811444325525444;
           Q=6.873773824708887; R=61.86396442237999; S=0.03050214047472335; T=0.03050214047472335 }
    9 = @{ I=0.6339497749751537; J=0.6339497749751537;
           O=0.01415294505639593; P=0.01415294505639593;
           S=0.008972256333737912; T=0.008972256333737914 }
    10 = @{ I=0.6339497749751537; J=0.6339497749751537;
            O=0.9377326116883496; P=0.9377326116883496;
            S=0.5944753781666925; T=0.5944753781666925 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
